# Revert "Powerpoint writer: consolidate text run nodes."
#
# The author's previous commit had merged runs like
#   <a:r><a:t>chicken</a:t></a:r><a:r><a:t> </a:t></a:r><a:r><a:t>and</a:t></a:r>...
# into
#   <a:r><a:t>chicken </a:t></a:r><a:r><a:t>and </a:t></a:r>...
# (folding the trailing space into the previous word's run). This
# reverts that: every run that ends with a space immediately before
# another run gets split back into a "word" run and a standalone
# " " run, leaving a run-per-token structure with the space as its
# own run. The last token of a paragraph remains attached to nothing
# (no trailing split needed since there's no subsequent run to merge
# with).
#
# We can't directly splice <a:r> nodes from the COM object model, but
# assigning `.Text` on a `Characters(start, length)` sub-range forces
# the writer to break runs at those exact character boundaries (while
# leaving the (empty) run properties alone, since we assign back the
# same text that was already there). So, for each affected paragraph,
# we walk forward through the known tokens and re-assign each token's
# own text to itself -- this creates a run boundary right after every
# token that is followed by more text.

function Split-Runs {
    param(
        $TextRange,
        [int]$StartOffset,
        [string[]]$Tokens
    )

    $idx = $StartOffset
    foreach ($tok in $Tokens) {
        $len = $tok.Length
        if ($len -gt 0) {
            $chars = $TextRange.Characters($idx, $len)
            $chars.Text = $tok
        }
        $idx += $len
    }
}

$p = $ppt.ActivePresentation

# --- Slide 1 notes: "chicken and dumplings" ---
$notes1 = $p.Slides.Item(1).NotesPage.Shapes.Item(2).TextFrame.TextRange
Split-Runs $notes1 1 @('chicken', ' ', 'and', ' ')

# --- Slide 2 notes: "foo bar" ---
$notes2 = $p.Slides.Item(2).NotesPage.Shapes.Item(2).TextFrame.TextRange
Split-Runs $notes2 1 @('foo', ' ')

# --- Slide 3 notes: two paragraphs ---
#   "Some notes inside a column"
#   (blank paragraph)
#   "Some notes outside the column"
$notes3 = $p.Slides.Item(3).NotesPage.Shapes.Item(2).TextFrame.TextRange
Split-Runs $notes3 1 @('Some', ' ', 'notes', ' ', 'inside', ' ', 'a', ' ')
# "Some notes inside a column" = 27 chars, then a paragraph mark (\r),
# then an empty paragraph + its mark (\r), then the second paragraph
# starts at offset 1 + 27 + 1 + 1 = 30.
Split-Runs $notes3 30 @('Some', ' ', 'notes', ' ', 'outside', ' ', 'the', ' ')

# --- Slide 1 (body): "The moon" ---
$title1 = $p.Slides.Item(1).Shapes.Item(2).TextFrame.TextRange
Split-Runs $title1 1 @('The', ' ')

# --- Slide 2 (body): "Demonstration of simple table syntax, with alignment" ---
$caption2 = $p.Slides.Item(2).Shapes.Item(2).TextFrame.TextRange
Split-Runs $caption2 1 @('Demonstration', ' ', 'of', ' ', 'simple', ' ', 'table', ' ', 'syntax,', ' ', 'with', ' ')
